# Regenerate handoff report: new GUID / new commit-hash file names and
# refreshed handoff timestamps for the zh-cn and de-de localization targets.

$wb = $excel.ActiveWorkbook

$oldGuidFile = "27bb69cd-bfbe-4d8b-9ac5-24d34b940c52.md"
$newGuidFile = "189a0fc2-e894-4626-9d39-3aa246d0bdad.md"

$oldZhXlf = "27bb69cd-bfbe-4d8b-9ac5-24d34b940c52.2807ab06bc834e4a58611285b59cf5be1af38cd1.zh-cn.xlf"
$newZhXlf = "189a0fc2-e894-4626-9d39-3aa246d0bdad.f6bc6adf31fbfea0b5018aa1abcae6bff66d0cdf.zh-cn.xlf"

$oldDeXlf = "27bb69cd-bfbe-4d8b-9ac5-24d34b940c52.2807ab06bc834e4a58611285b59cf5be1af38cd1.de-de.xlf"
$newDeXlf = "189a0fc2-e894-4626-9d39-3aa246d0bdad.f6bc6adf31fbfea0b5018aa1abcae6bff66d0cdf.de-de.xlf"

$oldZhDate = "2016-03-10 14:50:14"
$newZhDate = "2016-03-10 14:50:36"

$oldDeDate = "2016-03-10 14:50:19"
$newDeDate = "2016-03-10 14:50:41"

foreach ($ws in $wb.Worksheets) {

    # A2 always holds the handback markdown file name (e2e/<guid>.md), on
    # every sheet (Overview, zh-cn, de-de). (Note: ".Value" getter is not
    # reliable for comparisons in this runtime, so read back via ".Text".)
    if ($ws.Range("A2").Text -eq $oldGuidFile) {
        $ws.Range("A2").Value = $newGuidFile
    }

    # C2 holds the latest handoff xlf file name, only present on the
    # per-locale sheets.
    if ($ws.Range("C2").Text -eq $oldZhXlf) {
        $ws.Range("C2").Value = $newZhXlf
    }
    if ($ws.Range("C2").Text -eq $oldDeXlf) {
        $ws.Range("C2").Value = $newDeXlf
    }

    # D2 holds the latest handoff datetime, only present on the per-locale
    # sheets.
    if ($ws.Range("D2").Text -eq $oldZhDate) {
        $ws.Range("D2").Value = $newZhDate
    }
    if ($ws.Range("D2").Text -eq $oldDeDate) {
        $ws.Range("D2").Value = $newDeDate
    }

    # Keep the hyperlink "display" text in sync with the cell values above.
    foreach ($h in $ws.Hyperlinks) {
        if ($h.TextToDisplay -eq $oldGuidFile) {
            $h.TextToDisplay = $newGuidFile
        }
        elseif ($h.TextToDisplay -eq $oldZhXlf) {
            $h.TextToDisplay = $newZhXlf
        }
        elseif ($h.TextToDisplay -eq $oldDeXlf) {
            $h.TextToDisplay = $newDeXlf
        }
    }
}
